# Fruta / hortaliza, semanal
# A new weekly record is inserted at the top of the "Ajo" data block
# (row 20), and the previously existing rows 20-33 shift down into rows
# 21-34 (row 34's former data simply falls off the bottom, matching the
# source diff which leaves row 35+ untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for the columns that move: D (Fecha),
# J (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio
# ponderado), P (Precio $/Kg) for rows 20 through 33.
$cols = @("D", "J", "K", "L", "M", "P")

$snapshot = @{}
for ($r = 20; $r -le 33; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Shift rows 20-33 down into rows 21-34.
for ($r = 33; $r -ge 20; $r--) {
    foreach ($c in $cols) {
        $ws.Range("$c$($r + 1)").Value2 = $snapshot[$r][$c]
    }
}

# Write the brand-new record for row 20 (this week's entry).
$ws.Range("D20").Value2 = 45162
$ws.Range("J20").Value2 = 400
$ws.Range("K20").Value2 = 17000
$ws.Range("L20").Value2 = 18000
$ws.Range("M20").Value2 = 17500
$ws.Range("P20").Value2 = 1750
